$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

$ws.Range("A3").Select()
